$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -0.7157805844231913

$ws.Range("B4").Value = -0.7002365474293518
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.6913962355556218

$ws.Range("B8").Value = -0.7848563185571278
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 0
